# ---------------------------------------------------------------------------
# Applies the "Idk changes i forgot?" edit to Guion Texto.docx:
#   1. Inserts a new empty paragraph (lang=es-ES) just before the
#      "(Pirata): ¡Vaya! El mapa me..." list item.
#   2. Marks a (moved) rendered-page-break right before the
#      "(Pirata): ¡Así que la división..." line, and drops the _GoBack
#      bookmark into the middle of that paragraph's text (splitting the
#      run "...y como" / "o una fracción...").
#   3. Removes the old rendered-page-break that used to sit in front of
#      "Encontrar Primer Cofre".
#   4. Adds a rendered-page-break in front of "-(Pirata): He encontrado
#      una de las 4 gemas...".
#   5. The _GoBack bookmark move in step 2 automatically vacates its old
#      slot in front of "El pirata abrió la última cueva..." (Word only
#      ever keeps one _GoBack bookmark, so re-adding it elsewhere deletes
#      the old one).
# ---------------------------------------------------------------------------

function Insert-WordParagraphXml($range, [string]$innerXml) {
    # Wraps $innerXml (the contents of a <w:p>) in a minimal single-part
    # WordProcessingML package and inserts it at $range. Word/this runtime
    # treats a <w:p> that starts exactly at a paragraph's own start
    # position as "flowing into" that paragraph instead of forcing a new
    # paragraph break, which is what lets us splice extra leading runs
    # (e.g. a lone <w:lastRenderedPageBreak/>) onto an existing paragraph.
    $pkg = "<?xml version='1.0'?>" +
        "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
        "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
        "<pkg:xmlData>" +
        "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:body><w:p>" + $innerXml + "</w:p></w:body></w:document>" +
        "</pkg:xmlData></pkg:part></pkg:package>"
    $range.InsertXML($pkg)
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) New empty paragraph right before "(Pirata): ¡Vaya! El mapa me..."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("muestra un Presiona M para abrir el mapa se lee el siguiente texto", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insPoint = $d.Range($rng.End, $rng.End)
Insert-WordParagraphXml $insPoint "<w:pPr><w:rPr><w:lang w:val='es-ES'/></w:rPr></w:pPr>"

# ---------------------------------------------------------------------
# 2a) lastRenderedPageBreak in front of "(Pirata): ¡Así que la división..."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("sí que la división de esos números también puedo verlo", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $d.Range($rng.Start, $rng.Start)
$para.Expand(4) | Out-Null
$paraStart = $d.Range($para.Start, $para.Start)
Insert-WordParagraphXml $paraStart "<w:r><w:rPr><w:i/><w:lang w:val='es-ES'/></w:rPr><w:lastRenderedPageBreak/></w:r>"

# ---------------------------------------------------------------------
# 2b) Split "...y como"/"o una fraccion..." and drop the _GoBack bookmark
#     in between. Adding a bookmark mid-run naturally splits the run, and
#     re-adding "_GoBack" relocates the document's single _GoBack bookmark
#     here (vacating wherever it used to be -- see step 5 below, which
#     needs no separate action because of this).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("o una fracción! Eso sí que no me lo esperaba.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $d.Range($rng.Start, $rng.Start)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null

# ---------------------------------------------------------------------
# 3) Remove the lastRenderedPageBreak that used to precede
#    "Encontrar Primer Cofre" (it moved earlier in the doc, see 2a).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Encontrar Primer Cofre", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Delete()
$insPoint = $d.Range($rng.Start, $rng.Start)
Insert-WordParagraphXml $insPoint "<w:r><w:rPr><w:b/><w:lang w:val='es-ES'/></w:rPr><w:t>Encontrar Primer Cofre</w:t></w:r>"

# ---------------------------------------------------------------------
# 4) lastRenderedPageBreak in front of "-(Pirata): He encontrado una de
#    las 4 gemas..."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("He encontrado una de las 4 gemas para abrir la puerta del tesoro final", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $d.Range($rng.Start, $rng.Start)
$para.Expand(4) | Out-Null
$paraStart = $d.Range($para.Start, $para.Start)
Insert-WordParagraphXml $paraStart "<w:r><w:rPr><w:i/><w:lang w:val='es-ES'/></w:rPr><w:lastRenderedPageBreak/></w:r>"

Write-Output "edit complete"
